$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 34

$ws.Range("F8").Value = 5
$ws.Range("H8").Value = 5

$ws.Range("F24").Value = 10
$ws.Range("H24").Value = 10

$ws.Range("F25").Value = 6
$ws.Range("H25").Value = 6

$ws.Range("E28").Value = 10
$ws.Range("F28").Value = 6
$ws.Range("H28").Value = 6

$ws.Range("F36").Value = 24
$ws.Range("H36").Value = 24

$ws.Range("E38").Value = 49

$ws.Range("E45").Value = 19

$ws.Range("F46").Value = 6
$ws.Range("H46").Value = 6

$ws.Range("E49").Value = 52

$ws.Range("F63").Value = 5
$ws.Range("H63").Value = 5

$ws.Range("E65").Value = 24

$ws.Range("F66").Value = 14
$ws.Range("H66").Value = 14

$ws.Range("E67").Value = 29
$ws.Range("F67").Value = 18
$ws.Range("H67").Value = 18

$ws.Range("E70").Value = 26

$ws.Range("E71").Value = 22

$ws.Range("F72").Value = 15
$ws.Range("H72").Value = 15

$ws.Range("E76").Value = 31

$ws.Range("E77").Value = 35
$ws.Range("F77").Value = 13
$ws.Range("H77").Value = 13

$ws.Range("E79").Value = 19

$ws.Range("E80").Value = 18

$ws.Range("F81").Value = 4
$ws.Range("H81").Value = 4

$ws.Range("E82").Value = 6

$ws.Range("E88").Value = 12
